$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before A; this shifts the existing A:AC data to B:AD.
$ws.Columns("A:A").Insert()

# New column A holds a "Match ID" field: a header label in row 3 (the real
# header row) and a constant value (5) for every data row, including the
# hidden totals row 20. Rows 1-2 (the hidden pandas multi-index header rows)
# get no value in column A, matching the rest of that row's blank styling.
$ws.Range("A3").Value = "Match ID"
$ws.Range("A4:A20").Value = 5

# Row 3's header and the per-player rows (4-19) use a bold, borderless style
# (same font as the existing bold header style, but without the border/
# alignment that the column-header style carries).
$ws.Range("A3:A19").Font.Bold = $true

# The COM layer stamps a spurious custom row-height on row 20 the moment its
# existing (hidden-row) cell is rewritten; AutoFit puts the row back to the
# sheet's default auto height so the row tag stays clean.
$ws.Rows("20:20").AutoFit()

# Restore the selection to the new Match ID column per the saved view.
$ws.Range("A3:A19").Select()
